# altChunk: fix missing page break.
#
# The outer document gets a continuous section break right after the first
# paragraph ("Outer para 1"/"outer, before sect break"), so that the altChunk
# (inner sub-document) starts life in its own section that mimics the
# overall page geometry (A4, 708 twip margins/header/footer/cols, restart
# page numbering at 1). The trailing sectPr (covering the altChunk) becomes
# an explicit "continuous" section as well, and picks up the same page
# geometry.

$d = $word.ActiveDocument

# --- 1. Insert a continuous section break at the end of paragraph 1 -------
# Word always mints a fresh (empty) paragraph to hold the new sectPr, so
# immediately merge that empty paragraph back into paragraph 1 by deleting
# the paragraph mark between them -- this leaves the original paragraph's
# text intact while its own paragraph mark now carries the new <w:sectPr>.
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$r.Collapse(0)
$r.InsertBreak(2)

$p1 = $d.Paragraphs(1)
$p1r = $p1.Range
$mark = $d.Range($p1r.End - 1, $p1r.End)
$mark.Delete()

# --- 2. New (first) section: page geometry + restart page numbering -------
$sec1 = $d.Sections(1)

$hf1 = $sec1.Headers(1)
$hf1.PageNumbers.RestartNumberingAtSection = $true
$hf1.PageNumbers.StartingNumber = 1

$ps1 = $sec1.PageSetup
$ps1.PageWidth = 595.3
$ps1.PageHeight = 841.9
$ps1.TopMargin = 70.85
$ps1.BottomMargin = 70.85
$ps1.LeftMargin = 70.85
$ps1.RightMargin = 70.85
$ps1.HeaderDistance = 35.4
$ps1.FooterDistance = 35.4
$ps1.Gutter = 0
$ps1.TextColumns.Spacing = 35.4

# --- 3. Trailing section (altChunk's sectPr): mark continuous + same geometry
$sec2 = $d.Sections(2)
$ps2 = $sec2.PageSetup
$ps2.SectionStart = 0  # wdSectionContinuous -> <w:type w:val="continuous"/>
$ps2.PageWidth = 595.3
$ps2.PageHeight = 841.9
$ps2.TopMargin = 70.85
$ps2.BottomMargin = 70.85
$ps2.LeftMargin = 70.85
$ps2.RightMargin = 70.85
$ps2.HeaderDistance = 35.4
$ps2.FooterDistance = 35.4
$ps2.Gutter = 0
$ps2.TextColumns.Spacing = 35.4

# --- 4. Reword the two outer paragraphs ------------------------------------
$d.Content.Find.Execute("Outer para 1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "outer, before sect break", 2)
$d.Content.Find.Execute("Outer para 2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "outer, after sect break", 2)
